# Adding a Tester profile
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# users sheet: add a new "Tester" user in row 7
# ----------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("users")

$wsUsers.Cells.Item(7, 1).Value = "test@test.com"
$wsUsers.Hyperlinks.Add($wsUsers.Cells.Item(7, 1), "mailto:test@test.com")
# Reuse the same visual style as the other email hyperlink cells (A2/A3)
$wsUsers.Cells.Item(2, 1).Copy()
$wsUsers.Cells.Item(7, 1).PasteSpecial(-4122)

$wsUsers.Cells.Item(7, 2).Value = 12345678
$wsUsers.Cells.Item(7, 3).Value = "Tester"
$wsUsers.Cells.Item(7, 4).Value = "Testing"
$wsUsers.Cells.Item(7, 7).Value = $true
$wsUsers.Cells.Item(7, 8).Value = "admin"

$wsUsers.Activate()
$wsUsers.Range("L7").Select()

# ----------------------------------------------------------------------
# watchlist sheet: add 4 watchlist entries (user_id 6) in rows 18-21
# ----------------------------------------------------------------------
$wsWatchlist = $wb.Worksheets.Item("watchlist")

$wsWatchlist.Cells.Item(18, 1).Value = 6
$wsWatchlist.Cells.Item(18, 2).Value = "TSLA"

$wsWatchlist.Cells.Item(19, 1).Value = 6
$wsWatchlist.Cells.Item(19, 2).Value = "AAPL"

$wsWatchlist.Cells.Item(20, 1).Value = 6
$wsWatchlist.Cells.Item(20, 2).Value = "MSFT"

$wsWatchlist.Cells.Item(21, 1).Value = 6
$wsWatchlist.Cells.Item(21, 2).Value = "DELL"

$wsWatchlist.Activate()
$wsWatchlist.Range("B21").Select()

# ----------------------------------------------------------------------
# portfolio sheet: add Tester's portfolio in row 7
# ----------------------------------------------------------------------
$wsPortfolio = $wb.Worksheets.Item("portfolio")

$wsPortfolio.Cells.Item(7, 1).Value = 6
$wsPortfolio.Cells.Item(7, 2).Value = "Tester's Portfolio"
$wsPortfolio.Cells.Item(7, 3).Value = 1000000

$wsPortfolio.Activate()
$wsPortfolio.Range("C9").Select()

# ----------------------------------------------------------------------
# portfolioprice sheet: add a price snapshot in row 7
# ----------------------------------------------------------------------
$wsPortfolioPrice = $wb.Worksheets.Item("portfolioprice")

$wsPortfolioPrice.Cells.Item(7, 1).Value = 6
$wsPortfolioPrice.Cells.Item(7, 2).Value = 6
$wsPortfolioPrice.Cells.Item(7, 3).Value = 1000000
$wsPortfolioPrice.Cells.Item(7, 4).Value = 0
# Reuse the date/time style already applied to the column (row 6)
$wsPortfolioPrice.Cells.Item(6, 5).Copy()
$wsPortfolioPrice.Cells.Item(7, 5).PasteSpecial(-4122)
$wsPortfolioPrice.Cells.Item(7, 5).Value = 44120.275080960651

$wsPortfolioPrice.Activate()
$wsPortfolioPrice.Range("F7").Select()

# ----------------------------------------------------------------------
# transaction sheet: fill in rows 12-13 (previously blank placeholders)
# ----------------------------------------------------------------------
$wsTransaction = $wb.Worksheets.Item("transaction")

$wsTransaction.Cells.Item(12, 1).Value = 6
$wsTransaction.Cells.Item(12, 2).Value = 6
$wsTransaction.Cells.Item(12, 3).Value = "DELL"
$wsTransaction.Cells.Item(12, 4).Value = 60
$wsTransaction.Cells.Item(12, 5).Value = 44119.483414351853
$wsTransaction.Cells.Item(12, 6).Value = 2310
$wsTransaction.Cells.Item(12, 7).Value = 0

$wsTransaction.Cells.Item(13, 1).Value = 6
$wsTransaction.Cells.Item(13, 2).Value = 6
$wsTransaction.Cells.Item(13, 3).Value = "TSLA"
$wsTransaction.Cells.Item(13, 4).Value = 60
$wsTransaction.Cells.Item(13, 5).Value = 44119.483414351853
$wsTransaction.Cells.Item(13, 6).Value = 3420
$wsTransaction.Cells.Item(13, 7).Value = 0

$wsTransaction.Activate()
$wsTransaction.Range("D14").Select()

# ----------------------------------------------------------------------
# exchanges sheet: becomes the active/selected tab
# ----------------------------------------------------------------------
$wsExchanges = $wb.Worksheets.Item("exchanges")
$wsExchanges.Activate()
